$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new IPO subscription record (하나32호스팩) needs to be inserted as the
# newest entry at the top of the data table, pushing the existing rows
# (2-16) down to (3-17).
$ws.Rows("2:2").Insert()

# The Insert() copies the header row's bold/centered/bordered style onto
# the new row; reset it back to the plain "Normal" style used by every
# other data row so the new row matches its siblings.
$ws.Range("A2:T2").Style = "Normal"

# Dates in this sheet are stored as plain text (e.g. "2024-03-14"), not
# as real Excel date serials. Force the text format on the three date
# columns before assigning so Excel doesn't auto-convert them to dates.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"

$ws.Range("A2").Value = "2024-03-18"
$ws.Range("B2").Value = "하나32호스팩"
$ws.Range("C2").Value = "하나"
$ws.Range("D2").Value = "2024-03-21"
$ws.Range("E2").Value = "2024-03-27"
$ws.Range("F2").Value = 6000000
$ws.Range("G2").Value = 3000000
$ws.Range("H2").Value = "-"
$ws.Range("I2").Value = 2000
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = "-"
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = "-"
$ws.Range("N2").Value = "-"
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = "-"
$ws.Range("Q2").Value = "-"
$ws.Range("R2").Value = "2389.8 : 1"
$ws.Range("S2").Value = "-"
$ws.Range("T2").Value = "-"
